# Update coding scheme for consistency
# Applies text revisions to the coding-scheme explanations, resizes the
# affected rows to fit the revised (longer) text, and restores the
# previously-active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell text content ----------------------------------------
# NOTE: cells are written in the same order the original author touched
# them so that the rebuilt shared-string table lands in the same order.

$ws.Range("D54").Value = "Conduct a within-site search on Zenodo, Figshare, or Software Heritage using the software name plus possible keywords as search term. If an archival copy of the software (no matter which version) is found, then this code is TRUE. Using the link of identified working repository for searching inside Software Heritage is particularly helpful.`nAlso examine the web search results when searching for the specific piece of software using possible search terms. If an archival copy in an institutional repository or in locations mentioned above is located in the search engine results page, then this code is also TRUE."
$ws.Range("D2").Value = "A specific computational product that needs to be instantiated by code to realize the reported research in the publication. Examples include specific computational models or algorithms that are implemented according to the context. Exceptions include a general analysis method being mentioned but is not necessarily implemented by some form of code. If the extracted context does not mention software, then code all the rest codes as `"0`"."
$ws.Range("D29").Value = "Use available information in in-text mention and references to search for the software for available online records. If no version is mentioned in the extracted texts, then code C3 as FALSE."
$ws.Range("C28").Value = "The software has at least one findable official presence (e.g., source code, online manual, publication, or an online resource such as a metadata record or webpage that is dedicated to the software)."
$ws.Range("D26").Value = "Usually if the reference is a `"software publication`" or discusses the software substantially, the authors are counted as publishers/creators of the software."

# --- Resize rows whose wrapped text now needs more vertical space ----
$ws.Rows.Item(2).RowHeight = 145
$ws.Rows.Item(26).RowHeight = 43.5
$ws.Rows.Item(28).RowHeight = 72.5
$ws.Rows.Item(54).RowHeight = 188.5

# --- Restore selection / scroll position --------------------------------
$ws.Range("D27").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 23
$win.ScrollColumn = 1
